# Apply the "Trade #14 closed" update across the workbook.
$wb = $excel.ActiveWorkbook

# ---- Summary sheet ----
$wsSummary = $wb.Worksheets.Item("Summary")
$wsSummary.Range("B5").Value = 0.04    # Total P&L %
$wsSummary.Range("B6").Value = 14      # Total Trades
$wsSummary.Range("B9").Value = 42.86   # Win Rate %

# ---- Strategy Status sheet ----
$wsStatus = $wb.Worksheets.Item("Strategy Status")
$wsStatus.Range("D4").Value = 14       # Trades (MarketMaking)
$wsStatus.Range("G4").Value = 42.86    # Win Rate % (MarketMaking)

# ---- All Trades sheet: close out trade #14 (row 15) ----
$wsTrades = $wb.Worksheets.Item("All Trades")
$wsTrades.Range("G15").Value = 0.09          # Exit Price
$wsTrades.Range("H15").Value = "CLOSED"      # Status
$wsTrades.Range("K15").Value = 100.03        # Capital After
$wsTrades.Range("P15").Value = "early_exit"  # Exit Reason
$wsTrades.Range("Q15").Value = 0.13          # Duration (min)

# ---- MarketMaking sheet: mirror the same trade row ----
$wsMM = $wb.Worksheets.Item("MarketMaking")
$wsMM.Range("G15").Value = 0.09
$wsMM.Range("H15").Value = "CLOSED"
$wsMM.Range("K15").Value = 100.03
$wsMM.Range("P15").Value = "early_exit"
$wsMM.Range("Q15").Value = 0.13
